$d = $word.ActiveDocument

# Update the title/date paragraph
$d.Content.Find.Execute("2024-04-23 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-24 Wednesday", 2)

# Update the division problems in the table, addressed by explicit
# row/column so that duplicate values among old/new text never
# collide with each other during replacement.
$tbl = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $newText
}

Set-CellText $tbl 1 2 "94÷4="
Set-CellText $tbl 1 3 "57÷3="
Set-CellText $tbl 1 4 "25÷2="
Set-CellText $tbl 1 5 "15÷9="

Set-CellText $tbl 5 1 "58÷2="
Set-CellText $tbl 5 2 "34÷2="
Set-CellText $tbl 5 3 "24÷6="
Set-CellText $tbl 5 4 "68÷8="
Set-CellText $tbl 5 5 "63÷2="

Set-CellText $tbl 9 1 "75÷8="
Set-CellText $tbl 9 2 "87÷3="
Set-CellText $tbl 9 3 "41÷5="
Set-CellText $tbl 9 4 "12÷6="
Set-CellText $tbl 9 5 "62÷3="

Set-CellText $tbl 13 1 "10÷4="
Set-CellText $tbl 13 2 "59÷9="
Set-CellText $tbl 13 3 "18÷3="
Set-CellText $tbl 13 4 "93÷5="
Set-CellText $tbl 13 5 "27÷7="

Set-CellText $tbl 17 1 "72÷4="
Set-CellText $tbl 17 2 "60÷3="
Set-CellText $tbl 17 3 "98÷6="
Set-CellText $tbl 17 4 "66÷5="
Set-CellText $tbl 17 5 "57÷2="
